$d = $word.ActiveDocument

# --- Title paragraph: "...Immeuble ID: 4" -> "...Immeuble ID: 1" ---
# Replace only the trailing digit via a precise sub-range so that the
# straight apostrophe earlier in "l'Immeuble" is left completely untouched
# (Find/Execute's replacement text would otherwise get "smart-quoted").
$titlePara = $d.Paragraphs.Item(2)
$digitRange = $d.Range($titlePara.Range.End - 2, $titlePara.Range.End - 1)
$digitRange.Text = "1"

# --- Table row 2: ID Logement, Type Diagnostic, Date Diagnostic ---
$table = $d.Tables.Item(1)
$table.Cell(2, 1).Range.Text = "1"
$table.Cell(2, 2).Range.Text = "Amiante"
$table.Cell(2, 3).Range.Text = "2023-01-10"
